# Auto-generated edit script: update stock name table (rows 2-21) in HotStock_Top20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "平潭发展"
$ws.Range("B2").Value = "三花智控"
$ws.Range("C2").Value = "三花智控"
$ws.Range("A3").Value = "安泰科技"
$ws.Range("B3").Value = "安泰科技"
$ws.Range("C3").Value = "华建集团"
$ws.Range("A4").Value = "多氟多"
$ws.Range("B4").Value = "平潭发展"
$ws.Range("C4").Value = "平潭发展"
$ws.Range("A5").Value = "三花智控"
$ws.Range("B5").Value = "多氟多"
$ws.Range("C5").Value = "华胜天成"
$ws.Range("A6").Value = "海峡创新"
$ws.Range("B6").Value = "长城军工"
$ws.Range("C6").Value = "方正科技"
$ws.Range("A7").Value = "达华智能"
$ws.Range("B7").Value = "海峡创新"
$ws.Range("C7").Value = "安泰科技"
$ws.Range("A8").Value = "大为股份"
$ws.Range("B8").Value = "工业富联"
$ws.Range("C8").Value = "长城军工"
$ws.Range("A9").Value = "长城军工"
$ws.Range("B9").Value = "方正科技"
$ws.Range("C9").Value = "达华智能"
$ws.Range("A10").Value = "方正科技"
$ws.Range("B10").Value = "大为股份"
$ws.Range("C10").Value = "多氟多"
$ws.Range("A11").Value = "工业富联"
$ws.Range("B11").Value = "闻泰科技"
$ws.Range("C11").Value = "楚江新材"
$ws.Range("A12").Value = "锦富技术"
$ws.Range("B12").Value = "达华智能"
$ws.Range("C12").Value = "大为股份"
$ws.Range("A13").Value = "楚江新材"
$ws.Range("B13").Value = "永鼎股份"
$ws.Range("C13").Value = "锦富技术"
$ws.Range("A14").Value = "胜宏科技"
$ws.Range("B14").Value = "楚江新材"
$ws.Range("C14").Value = "统一股份"
$ws.Range("A15").Value = "统一股份"
$ws.Range("B15").Value = "胜宏科技"
$ws.Range("C15").Value = "卧龙电驱"
$ws.Range("A16").Value = "天际股份"
$ws.Range("B16").Value = "格尔软件"
$ws.Range("C16").Value = "神开股份"
$ws.Range("A17").Value = "神开股份"
$ws.Range("B17").Value = "合力泰"
$ws.Range("C17").Value = "盈新发展"
$ws.Range("A18").Value = "格尔软件"
$ws.Range("B18").Value = "东方财富"
$ws.Range("C18").Value = "世龙实业"
$ws.Range("A19").Value = "永鼎股份"
$ws.Range("B19").Value = "首开股份"
$ws.Range("C19").Value = "海峡创新"
$ws.Range("A20").Value = "盈新发展"
$ws.Range("B20").Value = "卧龙电驱"
$ws.Range("C20").Value = "青岛双星"
$ws.Range("A21").Value = "卧龙电驱"
$ws.Range("B21").Value = "天际股份"
$ws.Range("C21").Value = "大洋电机"
